$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column D, shifting existing column D (and beyond) to the right.
$ws.Range("D1:F1").EntireColumn.Insert()

# Set header row values for the newly inserted columns.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Fill data rows 2-11 for the new columns with "NA".
$ws.Range("D2:F11").Value = "NA"

$ws.UsedRange | Out-Null
